# Add a new row (row 71) of data to each of the 4 worksheets, mirroring the
# existing "time / length / ID / ..." log format already present in rows 2-70.

$wb = $excel.ActiveWorkbook

# Per-sheet data for the new row 71, in column order A..I
$rowsData = @{
    "MID_LFT_#1" = @{
        A = 45857.46179398148
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x48"
        E = "0x07"
        F = 400
        G = [double]"5.68631262647113e+23"
        H = 328
        I = 7
    }
    "MID_LFT_#2" = @{
        A = 45857.46179398148
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x48"
        E = "0x19"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 328
        I = 25
    }
    "MID_PLT_#1" = @{
        A = 45857.46179398148
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x64"
        E = "0x15"
        F = 110
        G = [double]"5.68631262647113e+23"
        H = 100
        I = 15
    }
    "MID_PLT_#2" = @{
        A = 45857.46179398148
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x79"
        E = "0x9"
        F = 130
        G = [double]"5.68631262647113e+23"
        H = 121
        I = 9
    }
}

foreach ($sheetName in $rowsData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $rowsData[$sheetName]
    $newRow = 71

    # Column A keeps the same number format (date/time) as the rows above it.
    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
